# Refresh market-price / profit columns (H:N) on each job sheet, as produced by
# the scheduled market-data runner. Plain value updates - no formulas involved.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 19: Unbreak My Heart | Roof Tile
$ws.Range("H19").Value = 50816.527
$ws.Range("I19").Value = 83227.45
$ws.Range("J19").Value = 6251.5
$ws.Range("K19").Value = 83227.45
$ws.Range("L19").Value = 6251.5
$ws.Range("M19").Value = -83052.45
$ws.Range("N19").Value = -6601.5
# Row 43: Growing Is Knowing | Growth Formula Gamma
$ws.Range("H43").Value = 11853.4
$ws.Range("J43").Value = 9163.727999999999
$ws.Range("L43").Value = 9163.727999999999
$ws.Range("N43").Value = -9301.727999999999
# Row 112: Making Ends Meet | Superior Spiritbond Potion
$ws.Range("H112").Value = 2427.3333
$ws.Range("J112").Value = 2427.3333
$ws.Range("L112").Value = 7281.999899999999
$ws.Range("N112").Value = -9497.999899999999
# Row 137: Cutting Edge of Culinary Quality | Magnesia Whetstone
$ws.Range("H137").Value = 5141.4688
$ws.Range("I137").Value = 1497.421
$ws.Range("K137").Value = 4492.263
$ws.Range("M137").Value = -1942.263
# Row 138: All-night Crafting | Cunning Craftsman's Tisane
$ws.Range("H138").Value = 3369.0278
$ws.Range("J138").Value = 3536.2104
$ws.Range("L138").Value = 10608.6312
$ws.Range("N138").Value = -20888.6312

$ws = $wb.Worksheets.Item("ARM")
# Row 2: Ain't Got No Ingots | Bronze Ingot
$ws.Range("H2").Value = 4786479
$ws.Range("I2").Value = 6495642
$ws.Range("J2").Value = 822
$ws.Range("K2").Value = 6495642
$ws.Range("L2").Value = 822
$ws.Range("M2").Value = -6495529
$ws.Range("N2").Value = -1048
# Row 45: Hollow Hallmarks | Mythril Ingot
$ws.Range("H45").Value = 1303.9048
$ws.Range("I45").Value = 1387.9286
$ws.Range("K45").Value = 1387.9286
$ws.Range("M45").Value = -1010.9286
# Row 102: Smells of Rich Tama-hagane | Tama-hagane Ingot
$ws.Range("H102").Value = 1510
$ws.Range("I102").Value = 1510
$ws.Range("K102").Value = 1510
$ws.Range("M102").Value = 112
# Row 116: No Scope | Titanbronze Ingot
$ws.Range("H116").Value = 4786479
$ws.Range("I116").Value = 6495642
$ws.Range("J116").Value = 822
$ws.Range("K116").Value = 6495642
$ws.Range("L116").Value = 822
$ws.Range("M116").Value = -6493348
$ws.Range("N116").Value = -5410
# Row 132: Don't Bore Me, Ore Me | Mountain Chromite Ingot
$ws.Range("H132").Value = 6693.619
$ws.Range("I132").Value = 4926.613
$ws.Range("J132").Value = 11673.363
$ws.Range("K132").Value = 14779.839
$ws.Range("L132").Value = 35020.089
$ws.Range("M132").Value = -12249.839
$ws.Range("N132").Value = -40080.089

$ws = $wb.Worksheets.Item("BSM")
# Row 3: Hells Bells | Bronze Ingot
$ws.Range("H3").Value = 4786479
$ws.Range("I3").Value = 6495642
$ws.Range("J3").Value = 822
$ws.Range("K3").Value = 6495642
$ws.Range("L3").Value = 822
$ws.Range("M3").Value = -6495528
$ws.Range("N3").Value = -1050
# Row 12: A Hit Job | Bronze Chaser Hammer
$ws.Range("H12").Value = 14675.2
$ws.Range("I12").Value = 0
$ws.Range("J12").Value = 14675.2
$ws.Range("K12").Value = 0
$ws.Range("L12").Value = 14675.2
$ws.Range("M12").ClearContents()
$ws.Range("N12").Value = -15011.2
# Row 15: Anutha Spatha | Bronze Spatha
$ws.Range("H15").Value = 0
$ws.Range("I15").Value = 0
$ws.Range("K15").Value = 0
$ws.Range("M15").ClearContents()
# Row 64: With Bearings Straight | Mythrite Nugget
$ws.Range("H64").Value = 125962.375
$ws.Range("I64").Value = 1333
$ws.Range("K64").Value = 1333
$ws.Range("M64").Value = -1108
# Row 67: Bearing the Brunt (L) | Mythrite Nugget
$ws.Range("H67").Value = 125962.375
$ws.Range("I67").Value = 1333
$ws.Range("K67").Value = 1333
$ws.Range("M67").Value = -553
# Row 134: Ruthenium Supremium | Ruthenium Ingot
$ws.Range("H134").Value = 6268.951
$ws.Range("I134").Value = 4796.778
$ws.Range("K134").Value = 14390.334
$ws.Range("M134").Value = -11855.334

$ws = $wb.Worksheets.Item("CRP")
# Row 58: You Do the Heavy Lifting | Mahogany Lumber
$ws.Range("H58").Value = 8764.076999999999
$ws.Range("I58").Value = 9996
$ws.Range("J58").Value = 8540.091
$ws.Range("K58").Value = 9996
$ws.Range("L58").Value = 8540.091
$ws.Range("M58").Value = -9793
$ws.Range("N58").Value = -8946.091
# Row 62: Splinter in the Sewers | Cedar Lumber
$ws.Range("H62").Value = 14933.533
$ws.Range("J62").Value = 17213.2
$ws.Range("L62").Value = 17213.2
$ws.Range("N62").Value = -18461.2
# Row 65: The Lumber of Their Discontent (L) | Cedar Lumber
$ws.Range("H65").Value = 14933.533
$ws.Range("J65").Value = 17213.2
$ws.Range("L65").Value = 86066
$ws.Range("N65").Value = -92306
# Row 136: Turali Quality | Dark Mahogany Lumber
$ws.Range("H136").Value = 8764.076999999999
$ws.Range("I136").Value = 9996
$ws.Range("J136").Value = 8540.091
$ws.Range("K136").Value = 29988
$ws.Range("L136").Value = 25620.273
$ws.Range("M136").Value = -27438
$ws.Range("N136").Value = -30720.273
# Row 141: No Greater Treasure | Claro Walnut Necklace of Gathering
$ws.Range("H141").Value = 186754.89
$ws.Range("J141").Value = 220029.9
$ws.Range("L141").Value = 220029.9
$ws.Range("N141").Value = -230389.9

$ws = $wb.Worksheets.Item("CUL")
# Row 4: In Hot Water | Boiled Egg
$ws.Range("H4").Value = 6739080.5
$ws.Range("I4").Value = 6500814
$ws.Range("J4").Value = 6977347
$ws.Range("K4").Value = 19502442
$ws.Range("L4").Value = 20932041
$ws.Range("M4").Value = -19502330
$ws.Range("N4").Value = -20932265
# Row 64: The Aroma of Faith | Baked Onion Soup
$ws.Range("H64").Value = 4896.6
$ws.Range("J64").Value = 4919.5
$ws.Range("L64").Value = 14758.5
$ws.Range("N64").Value = -15298.5
# Row 67: Soup's On (L) | Baked Onion Soup
$ws.Range("H67").Value = 4896.6
$ws.Range("J67").Value = 4919.5
$ws.Range("L67").Value = 14758.5
$ws.Range("N67").Value = -16630.5
# Row 137: Creative Chocolate | Gateau au Chocolat
$ws.Range("H137").Value = 10343.167
$ws.Range("I137").Value = 13460.625
$ws.Range("J137").Value = 4108.25
$ws.Range("K137").Value = 40381.875
$ws.Range("L137").Value = 12324.75
$ws.Range("M137").Value = -35281.875
$ws.Range("N137").Value = -22524.75

$ws = $wb.Worksheets.Item("GSM")
# Row 47: Wear Your Patriotic Pin | Peridot Choker
$ws.Range("H47").Value = 23505.357
$ws.Range("J47").Value = 23390.54
$ws.Range("L47").Value = 23390.54
$ws.Range("N47").Value = -24526.54
# Row 80: Needs More Prayerbell | Hardsilver Ingot
$ws.Range("H80").Value = 8338.333000000001
$ws.Range("I80").Value = 6137.5
$ws.Range("J80").Value = 10099
$ws.Range("K80").Value = 6137.5
$ws.Range("L80").Value = 10099
$ws.Range("M80").Value = -5139.5
$ws.Range("N80").Value = -12095
# Row 83: With a Noise That Reaches Heaven (L) | Hardsilver Ingot
$ws.Range("H83").Value = 8338.333000000001
$ws.Range("I83").Value = 6137.5
$ws.Range("J83").Value = 10099
$ws.Range("K83").Value = 30687.5
$ws.Range("L83").Value = 50495
$ws.Range("M83").Value = -25695.5
$ws.Range("N83").Value = -60479
# Row 126: Gold Rush Order | Phrygian Gold Ingot
$ws.Range("H126").Value = 4625.846
$ws.Range("I126").Value = 4664.364
$ws.Range("K126").Value = 13993.092
$ws.Range("M126").Value = -11523.092

$ws = $wb.Worksheets.Item("LTW")
# Row 25: A Rush on Ringbands | Hard Leather Ringbands
$ws.Range("H25").Value = 9917.833000000001
$ws.Range("I25").Value = 5669
$ws.Range("K25").Value = 5669
$ws.Range("M25").Value = -5439
# Row 61: Spelling Me Softly | Raptor Leather
$ws.Range("H61").Value = 3377.2273
$ws.Range("I61").Value = 3595.1
$ws.Range("J61").Value = 1198.5
$ws.Range("K61").Value = 3595.1
$ws.Range("L61").Value = 1198.5
$ws.Range("M61").Value = -3393.1
$ws.Range("N61").Value = -1602.5
# Row 113: Peace in Rest | Atrociraptor Leather
$ws.Range("H113").Value = 3377.2273
$ws.Range("I113").Value = 3595.1
$ws.Range("J113").Value = 1198.5
$ws.Range("K113").Value = 3595.1
$ws.Range("L113").Value = 1198.5
$ws.Range("M113").Value = -1425.1
$ws.Range("N113").Value = -5538.5

$ws = $wb.Worksheets.Item("WVR")
# Row 132: Comfy Cabins | Snow Cotton Cloth
$ws.Range("H132").Value = 4823.14
$ws.Range("I132").Value = 3905.3142
$ws.Range("K132").Value = 11715.9426
$ws.Range("M132").Value = -9185.942599999998
# Row 136: Weaving the Envelope | Sarcenet Cloth
$ws.Range("H136").Value = 2619.7144
$ws.Range("I136").Value = 1941.3636
$ws.Range("K136").Value = 5824.0908
$ws.Range("M136").Value = -3274.0908
